$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Delete the three device columns that are no longer in service
# (APPLE_iPhone7plus_iOS_13.5.1_b1cc7, APPLE_iPhoneXR_iOS_15.2.1_b3558
# and APPLE_iPhoneXS_iOS_14.3.0_33b29). Remove from right to left so
# earlier deletions don't shift the column letters of the ones still
# queued for removal.
$ws.Range("K1").EntireColumn.Delete()
$ws.Range("J1").EntireColumn.Delete()
$ws.Range("F1").EntireColumn.Delete()

# Re-anchor the "begins with Working / Not Working" conditional
# formatting to the new (narrower) row so it still spans the header
# row of device columns without losing its styling (dxf/priority).
$oldHeaderRow = $ws.Range("B2:M2")
$newHeaderRow = $ws.Range("B2:J2")
for ($i = 1; $i -le $oldHeaderRow.FormatConditions.Count; $i++) {
    $oldHeaderRow.FormatConditions.Item($i).ModifyAppliesToRange($newHeaderRow)
}

$ws.Activate()
[void]$ws.Range("G16").Select()
